# LocationFacetMapping.xlsx — apply location changes to facet
# Adds two new rows to the "LocationFacetMapping" sheet:
#   - "Mann Special Collections (Request in advance)"  -> "Mann Library > Special Collections"
#   - "Music Library (Lincoln Hall)"                    -> "Music Library > Main Collection"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LocationFacetMapping")

# --- Insert both new rows first (so row numbers below are easy to reason
# about), then populate their values. The Music Library shared string is
# authored before the Mann Special Collections one so the new entries land
# in the sharedStrings table in that same order.

# Row right after row 53 ("Mann Library Special Collections (Non-Circulating)"),
# pushing "Mann Serials" and everything below it down by one row.
$ws.Rows.Item(54).Insert()

# Row right before the (now shifted) "Cox Library of Music (Lincoln Hall)"
# row, which after the first insertion sits at row 62 -- so the new row
# goes in at row 61.
$ws.Rows.Item(61).Insert()

$ws.Range("A61").Value = "Music Library (Lincoln Hall)"
$ws.Range("E61").Value = "Music Library > Main Collection"

$ws.Range("A54").Value = "Mann Special Collections (Request in advance)"
$ws.Range("E54").Value = "Mann Library > Special Collections"

# --- Restore the view state (best effort)
$ws.Range("A48").Select() | Out-Null
